# Apply "average with safety stocks" edit:
#  - Productdata sheet: scale InventoryCosts (D), BackorderCosts (F) and
#    LostSale (I) columns for rows 2-11 by a factor of 0.0004
#  - ForcastedStandardDeviation sheet: zero out the AverageDemand /
#    StandardDevDemands columns (B:E) for rows 9-11

$wb = $excel.ActiveWorkbook

# --- Productdata sheet -----------------------------------------------
$wsProd = $wb.Worksheets.Item("Productdata")

$newD = @(0.0016, 0.0028, 0.0024, 0.0012, 0.0012, 0.0012, 0.0008, 0.0004, 0.0004, 0.0004)
$newF = @(0.0032, 0.0056, 0.0048, 0.0024, 0.0024, 0.0024, 0.0016, 0.0008, 0.0008, 0.0008)
$newI = @(0.032,  0.056,  0.048,  0.024,  0.024,  0.024,  0.016,  0.008,  0.008,  0.008)

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $wsProd.Range("D$row").Value = $newD[$i]
    $wsProd.Range("F$row").Value = $newF[$i]
    $wsProd.Range("I$row").Value = $newI[$i]
}

# --- ForcastedStandardDeviation sheet ---------------------------------
$wsStd = $wb.Worksheets.Item("ForcastedStandardDeviation")

foreach ($row in 9..11) {
    foreach ($col in @("B", "C", "D", "E")) {
        $wsStd.Range("$col$row").Value = 0
    }
}
